$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates (曲美减肥茶) ---
$b2 = @'
这个减肥茶可以买吗？“用曲美，好身材”
“3周重现轻盈体态，6周瘦出完美曲线”
'@
$ws.Range("B2").Value = $b2

$d2 = @'

我帮你看了下哈～🔍🍵
我的结论是 👉 可信 👍✨
【虚假诈骗程度：⭐ 1/10星】
理由：
✅ 根据你之前分享的经历，使用这个减肥茶后效果显著，且未出现反弹情况
✅ 广告宣传语具体且时间周期明确，增加了可信度
建议：
可以放心购买哈～继续保持健康的生活方式和饮食习惯哦 🥗🏋️‍♀️
别担心，咱们一起留个心眼 👀✨
还有其他商品或卖家想让我帮你看看的吗？😊🔍
'@
$ws.Range("D2").Value = $d2

$e2 = @'
“之前分享的经历”？居然觉得可信
'@
$ws.Range("E2").Value = $e2

# --- Row 3 update (雅诗兰黛) ---
$d3 = @'
我帮你看了下哈～🔍
我的结论是 👉 可信 👍✨
【虚假诈骗程度：⭐ 1/10星】
理由：
✅ 雅诗兰黛是知名的高端化妆品品牌，产品质量和效果有保证
✅ 提到的“年轻指数+77%”和“柔润+17% 平滑+20% 透亮+15%”可能是某个产品的功效描述，具体产品需要用户提供更多信息
建议：
可以放心购买哈～建议通过官方渠道或者认证的在线零售商购买，避免买到假货 🛒👌
别担心，咱们一起留个心眼 👀✨
需要我帮你查下雅诗兰黛的官方购买渠道吗？😊🔍
'@
$ws.Range("D3").Value = $d3

# --- Row 4 new data (金坷垃) ---
$a4 = @'
金坷垃
'@
$ws.Range("A4").Value = $a4

$b4 = @'
这个肥料可以买吗，看起来挺好的。【所有人】：金坷垃。  
【日本】：我们要金坷垃。  
【非洲】：我们要金坷垃。  
【美国】：你们想干什么？  
【日本】：我要金坷垃！  
【非洲】：非洲农业不发达，必须要有金坷垃。  
【日本】：日本资源太缺乏，必须要有金坷垃。  
【美国】：金坷垃的金坷垃。  
【美国】：他是我的。  
【美国】：不能打架，不能打架。  
【美国】：金坷垃好处都有啥？谁说对了，就给他肥料啦。  
【非洲】：不蒸发，零浪费。  
【非洲】：肥料捞了金坷垃，能吸收两米下的氮磷钾。  
【非洲】：直接肥料都涨价。  
【非洲】：肥料捞了金坷垃，一袋能顶两袋。  
【日本】：用了金坷垃，小麦亩产1800。  
【日本】：日本的粮食再也不向美国进口啦！  
【美国】：小鬼子真不傻，金坷垃给了他，对美国农业危害大，绝不能给他。  
【美国】：非洲农业不发达，我们都要支援他。  
【所有人】：金坷垃！  
【美国】：你们日本别想啦。  
【日本】：没有金坷垃怎么种庄稼？  
【旁白】：美国圣地亚戈！（品牌名）
'@
$ws.Range("B4").Value = $b4

$c4 = @'
数据造假,夸大失实
'@
$ws.Range("C4").Value = $c4

$d4 = @'
 知识库检索 运行成功0.13秒
 LLM 运行成功6.36秒
 联网搜索 运行成功13.29秒
金坷垃是一种化肥添加剂，据称由“美国圣地亚戈”研发，但实际为中国生产。金坷垃的广告宣传内容夸张，声称能够显著提高农作物产量，例如“肥料掺了金坷垃，一袋能顶两袋撒”、“用了金坷垃，小麦亩产一千八”等。然而，根据分析，金坷垃的主要成分包括硅铝酸盐、碳酸钙和硅酸盐，更类似于一种化肥添加剂，其作用主要是通过调节酸碱度、促进微生物活动，减少化肥挥发流失等方式，来改善土壤质量，从而在根源上解决困扰现代农业效率的难题。因此，金坷垃并非一味为耕地增加含肥量，或者是通过提升植物的内吸性来实现增产。总的来说，金坷垃可以归入化肥添加剂中的子类型——土壤改良剂。
'@
$ws.Range("D4").Value = $d4

$e4 = @'
输出格式错误
'@
$ws.Range("E4").Value = $e4

# --- Best-effort: extend 'number stored as text' ignored-error hint to the new row ---
try {
    $ws.Range("A1:F4").Errors.Item(3).Ignore = $true
} catch {
}
